# "Generate report" population: fill in the DOB (Date of Birth) column (C)
# for the employee rows with literal text values, matching the values
# produced by the report-generation feature.
#
# NumberFormat is temporarily set to Text ("@") before assignment so that
# Excel stores the values as literal strings (e.g. "01/01/2000") rather
# than auto-converting them into date serial numbers. ClearFormats() is
# then used to drop the temporary text-format style again, leaving the
# cells with their original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dobValues = @{
    "C2" = "01/01/2000"
    "C3" = "01/01/1980"
    "C4" = "01/01/1980"
}

foreach ($addr in $dobValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dobValues[$addr]
    $cell.ClearFormats()
}
